# NATMI TPM re-run: Tgfb2-Acvr1 ligand/receptor stats were recomputed with
# updated TPM input (one more cell now counts as expressing Tgfb2 in the
# "ECs" sending cluster: 2/3 -> 3/3 detection), which ripples through the
# per-cluster ligand totals/specificities (columns E-J), the shared
# per-cluster receptor totals/specificities (columns M-P), and the derived
# edge-weight columns (Q-T) for every row that references those clusters.
# Values below are the recomputed outputs written back into the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.609586333333333
$ws.Range("H2").Value = 4.828759
$ws.Range("I2").Value = 0.05107819292772156
$ws.Range("J2").Value = 0.05107819292772156
$ws.Range("M2").Value = 4.621579
$ws.Range("N2").Value = 13.864737
$ws.Range("O2").Value = 0.1778708528171788
$ws.Range("P2").Value = 0.1778708528171788
$ws.Range("Q2").Value = 7.438830396820332
$ws.Range("R2").Value = 66.949473571383
$ws.Range("S2").Value = 0.009085321736414225
$ws.Range("T2").Value = 0.009085321736414225
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.609586333333333
$ws.Range("H3").Value = 4.828759
$ws.Range("I3").Value = 0.05107819292772156
$ws.Range("J3").Value = 0.05107819292772156
$ws.Range("N3").Value = 46.543441
$ws.Range("O3").Value = 0.5971062807549863
$ws.Range("P3").Value = 0.5971062807549863
$ws.Range("Q3").Value = 24.97189551330211
$ws.Range("R3").Value = 224.747059619719
$ws.Range("S3").Value = 0.03049910980675747
$ws.Range("T3").Value = 0.03049910980675747
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.609586333333333
$ws.Range("H4").Value = 4.828759
$ws.Range("I4").Value = 0.05107819292772156
$ws.Range("J4").Value = 0.05107819292772156
$ws.Range("O4").Value = 0.2250228664278349
$ws.Range("P4").Value = 0.2250228664278349
$ws.Range("Q4").Value = 9.41079953376911
$ws.Range("R4").Value = 84.69719580392199
$ws.Range("S4").Value = 0.01149376138454987
$ws.Range("T4").Value = 0.01149376138454987
$ws.Range("I5").Value = 0.5992082897496871
$ws.Range("J5").Value = 0.5992082897496871
$ws.Range("M5").Value = 4.621579
$ws.Range("N5").Value = 13.864737
$ws.Range("O5").Value = 0.1778708528171788
$ws.Range("P5").Value = 0.1778708528171788
$ws.Range("Q5").Value = 87.26637698644068
$ws.Range("R5").Value = 785.3973928779661
$ws.Range("S5").Value = 0.1065816895129
$ws.Range("T5").Value = 0.1065816895129
$ws.Range("I6").Value = 0.5992082897496871
$ws.Range("J6").Value = 0.5992082897496871
$ws.Range("N6").Value = 46.543441
$ws.Range("O6").Value = 0.5971062807549863
$ws.Range("P6").Value = 0.5971062807549863
$ws.Range("S6").Value = 0.3577910332899918
$ws.Range("T6").Value = 0.3577910332899918
$ws.Range("I7").Value = 0.5992082897496871
$ws.Range("J7").Value = 0.5992082897496871
$ws.Range("O7").Value = 0.2250228664278349
$ws.Range("P7").Value = 0.2250228664278349
$ws.Range("S7").Value = 0.1348355669467953
$ws.Range("T7").Value = 0.1348355669467953
$ws.Range("H8").Value = 33.060729
$ws.Range("I8").Value = 0.3497135173225914
$ws.Range("J8").Value = 0.3497135173225914
$ws.Range("M8").Value = 4.621579
$ws.Range("N8").Value = 13.864737
$ws.Range("O8").Value = 0.1778708528171788
$ws.Range("P8").Value = 0.1778708528171788
$ws.Range("Q8").Value = 50.930923623697
$ws.Range("R8").Value = 458.378312613273
$ws.Range("S8").Value = 0.06220384156786456
$ws.Range("T8").Value = 0.06220384156786458
$ws.Range("H9").Value = 33.060729
$ws.Range("I9").Value = 0.3497135173225914
$ws.Range("J9").Value = 0.3497135173225914
$ws.Range("N9").Value = 46.543441
$ws.Range("O9").Value = 0.5971062807549863
$ws.Range("P9").Value = 0.5971062807549863
$ws.Range("Q9").Value = 170.9733432920544
$ws.Range("S9").Value = 0.208816137658237
$ws.Range("T9").Value = 0.208816137658237
$ws.Range("H10").Value = 33.060729
$ws.Range("I10").Value = 0.3497135173225914
$ws.Range("J10").Value = 0.3497135173225914
$ws.Range("O10").Value = 0.2250228664278349
$ws.Range("P10").Value = 0.2250228664278349
$ws.Range("Q10").Value = 64.43226780613134
$ws.Range("R10").Value = 579.890410255182
$ws.Range("S10").Value = 0.07869353809648984
$ws.Range("T10").Value = 0.07869353809648984
